# Increase font sizes across the resume document.
# Mapping is based on the point size (not half-points) of each run:
#   16 -> 18   (name header)
#    9 -> 10   (contact info / body text / bullets / dates-locations)
#   12 -> 13   (section headers)
#   11 -> 12   (job titles)
#   10 -> 11   (overview / section intro paragraphs)
# Every paragraph in this document uses a single, uniform font size across
# all of its runs, so it is safe (and simplest) to resize at the paragraph
# Range level.

$d = $word.ActiveDocument

$sizeMap = @{
    16 = 18
    9  = 10
    12 = 13
    11 = 12
    10 = 11
}

foreach ($p in $d.Paragraphs) {
    # Exclude the trailing paragraph-mark character from the range so that
    # we only touch the actual run(s) of text, not the paragraph mark's
    # own run properties (which live in w:pPr/w:rPr).
    $start = $p.Range.Start
    $end = $p.Range.End - 1
    if ($end -le $start) {
        continue
    }
    $rng = $d.Range($start, $end)
    $current = $rng.Font.Size
    if ($sizeMap.ContainsKey($current)) {
        $rng.Font.Size = $sizeMap[$current]
    }
}
